# Regen sval data to filter save games
# Update the numeric columns B:G for rows 2-8 with the new computed values.
# Column A (dates) and row 1 (headers) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059)
    3 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455)
    4 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 0, 15.54439750930335)
    5 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
    6 = @(0.1169995834814548, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 0, 1.67637130870356)
    7 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 0, 19.48425592650926)
    8 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
